# Update iserv_stats.xlsx row 23 (2025-10) with refreshed figures.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Raw measured values for the month
$ws.Range("B23").Value = 6336
$ws.Range("C23").Value = 998
$ws.Range("D23").Value = 5899728

# Derived metrics: users per school, and year-over-year % changes vs. row 11 (12 months earlier)
$b23 = $ws.Range("B23").Value2
$c23 = $ws.Range("C23").Value2
$d23 = $ws.Range("D23").Value2
$b11 = $ws.Range("B11").Value2
$c11 = $ws.Range("C11").Value2
$d11 = $ws.Range("D11").Value2

$ws.Range("E23").Value = $d23 / $b23
$ws.Range("F23").Value = ($b23 - $b11) / $b11 * 100
$ws.Range("G23").Value = ($c23 - $c11) / $c11 * 100
$ws.Range("H23").Value = ($d23 - $d11) / $d11 * 100
